$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 11 (years 2000-2009), shifting rows 12-14 (2010-2012) up to 2-4
$ws.Rows("2:11").Delete()
